$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "IL6_IL6R"
$ws.Range("B2").Value = "IL6"
$ws.Range("C2").Value = "IL6R"
$ws.Range("D2").Value = "CCIDB_Human_new.xlsx; Human-2014-Qiao-LR-pairs.xlsx; Human-2015-Choi-LR-pairs.txt; Human-2015-Ramilowski-LR-pairs.xlsx; Human-2017-Pavlicev-LR-pairs.xlsx; Human-2019-Wang-LR-pairs.csv; Human-2019-Ximerakis-BaderLab-2017_new.xlsx; Human-2020-Hou-LR-pairs.xlsx; Human-2020-Shao-LR-pairs.txt; ICELLNETdb_v2_new.xlsx; LRdb_122019_SingleCellSignalIR.xlsx; NicheNet-LR-pairs.csv; OmniPathPPIs_new.tsv"
$ws.Range("E2").Value = 41
